$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46041
$ws.Range("B2").Value = 11794.6060739326
$ws.Range("C2").Value = 11287.1731960918
$ws.Range("D2").Value = 19275.86
$ws.Range("E2").Value = 7274.96806962966
$ws.Range("F2").Value = -29.738280594941

$ws.Range("A3").Value = 46042
$ws.Range("B3").Value = 12262.9170947069
$ws.Range("C3").Value = 12103.8365585735
$ws.Range("D3").Value = 12075.86
$ws.Range("E3").Value = 8292.30168747829
$ws.Range("F3").Value = 346.678260252158

$ws.Range("A4").Value = 46043
$ws.Range("B4").Value = 12847.6679828473
$ws.Range("C4").Value = 12654.5979892251
$ws.Range("D4").Value = 12075.86
$ws.Range("E4").Value = 8726.46425634832
$ws.Range("F4").Value = 387.716760232225

$ws.Range("A5").Value = 46044
$ws.Range("B5").Value = 13009.9244070811
$ws.Range("C5").Value = 12607.0743974187
$ws.Range("D5").Value = 12075.86
$ws.Range("E5").Value = 8848.16481300129
$ws.Range("F5").Value = 390.807467100833

$ws.Range("A6").Value = 46045
$ws.Range("B6").Value = 12876.9307050575
$ws.Range("C6").Value = 11708.0831074302
$ws.Range("D6").Value = 12075.86
$ws.Range("E6").Value = 8716.14642349942
$ws.Range("F6").Value = 347.848730455399

$ws.Range("A7").Value = 46046
$ws.Range("B7").Value = 5107.54318705847
$ws.Range("C7").Value = 8007.25018305176
$ws.Range("D7").Value = 12075.86
$ws.Range("E7").Value = 8264.44717782362
$ws.Range("F7").Value = 174.826556703141

$ws.Range("A8").Value = 46047
$ws.Range("B8").Value = 5000.01932310789
$ws.Range("C8").Value = 8102.26701796858
$ws.Range("D8").Value = 12075.86
$ws.Range("E8").Value = 8256.70262722195
$ws.Range("F8").Value = 178.462901882939

$ws.Range("A9").Value = 46048
$ws.Range("B9").Value = 12315.9682835607
$ws.Range("C9").Value = 11976.3433837723
$ws.Range("D9").Value = 12075.86
$ws.Range("E9").Value = 8314.28494857018
$ws.Range("F9").Value = 342.282013847605

$ws.Range("A10").Value = 46049
$ws.Range("B10").Value = 12315.9682835607
$ws.Range("C10").Value = 12384.7629330449
$ws.Range("D10").Value = 12075.86
$ws.Range("E10").Value = 8314.28494857018
$ws.Range("F10").Value = 359.299495067297

$ws.Range("A11").Value = 46050
$ws.Range("B11").Value = 12315.9682835607
$ws.Range("C11").Value = 12507.3824691386
$ws.Range("D11").Value = 12075.86
$ws.Range("E11").Value = 8314.28494857018
$ws.Range("F11").Value = 364.408642404531

$ws.Range("A12").Value = 46051
$ws.Range("B12").Value = 12315.9682835607
$ws.Range("C12").Value = 12519.3512424834
$ws.Range("D12").Value = 12075.86
$ws.Range("E12").Value = 8314.28494857018
$ws.Range("F12").Value = 364.907341293898

$ws.Range("A13").Value = 46052
$ws.Range("B13").Value = 12315.9682835607
$ws.Range("C13").Value = 11862.3067826234
$ws.Range("D13").Value = 12075.86
$ws.Range("E13").Value = 8314.28494857018
$ws.Range("F13").Value = 337.530488799732

$ws.Range("A14").Value = 46053
$ws.Range("B14").Value = 4867.38022112383
$ws.Range("C14").Value = 8854.09457363252
$ws.Range("D14").Value = 12075.86
$ws.Range("E14").Value = 7930.22316351309
$ws.Range("F14").Value = 196.185739047734

$ws.Range("A15").Value = 46054
$ws.Range("B15").Value = 5046.61366744637
$ws.Range("C15").Value = 8736.84096535956
$ws.Range("D15").Value = 9743.86
$ws.Range("E15").Value = 7890.29135676355
$ws.Range("F15").Value = 286.803013421796
